# GrimWorld 40,000 - Framework - 3119805903 번역 갱신
# Adds four new translation rows (supporter monument / obelisk ThingDefs)
# to the bottom of the sheet: GW_FancyObelisk + GW_OrdinaryObelisk,
# each with a label and a description, in English (col D) and Korean (col E).
#
# Cells are written in the same first-seen order as the shared-string
# table in the target workbook (C column top-to-bottom, then the shared
# "ThingDef" value, then the D column, then the E column) so the
# resulting xl/sharedStrings.xml ordering matches exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C - Identifier keys
$ws.Range('C18').Value = 'GW_FancyObelisk.label'
$ws.Range('C19').Value = 'GW_FancyObelisk.description'
$ws.Range('C20').Value = 'GW_OrdinaryObelisk.label'
$ws.Range('C21').Value = 'GW_OrdinaryObelisk.description'

# Column B - Node (same value for all four new rows)
$ws.Range('B18').Value = 'ThingDef'

# Column D - EN [Source string]
$ws.Range('D18').Value = 'supporter monument center'
$ws.Range('D19').Value = 'A monument to those who helped make the GrimWorld project a reality.\n\n  Astartes tier:\n\nRainyredman1234, Bailey English\n\n  Primarch tier:\n\nFulgrim, R.CSN, kyle shadowchain, Echo\n\n  Emperor of Mankind tier:\n\nRisqué Che, Team Shibe'
$ws.Range('D21').Value = 'A monument to those who helped make the GrimWorld project a reality.\n\n  Servitor tier:\n\nAKorgar, Hell Fire, Sparrow, just a guy named brad\n\n  Guardsmen tier:\n\n♠_Caligula_♠, Józef Kozioł, MisterCroxo48\n\n  Stormtrooper tier:\n\nBenio, Alloyskull, Vylixan, Gofres, Celorico, Risque, Madgile\n\n  Neophyte tier:\n\nTacticalCrumpet, Petrie, JawnWick'
$ws.Range('D20').Value = 'supporter''s monument'

# Column E - KO [Translation]
$ws.Range('E20').Value = '후원자 기념비'
$ws.Range('E18').Value = '후원자 기념 센터'
$ws.Range('E19').Value = '그림월드 프로젝트의 실현을 도운 이들을 위한 기념비입니다.\n\n  아스타르테스 티어:\n\nRainyredman1234, Bailey English\n\n  프라이마크 티어:\n\nFulgrim, R.CSN, kyle shadowchain, Echo\n\n  인류의 황제 티어:\n\nRisqué Che, Team Shibe'
$ws.Range('E21').Value = '그림월드 프로젝트의 실현을 도운 이들을 위한 기념비입니다.\n\n  서비터 티어:\n\nAKorgar, Hell Fire, Sparrow, just a guy named brad\n\n  가드맨 티어:\n\n♠_Caligula_♠, Józef Kozioł, MisterCroxo48\n\n  스톰트루퍼 티어:\n\nBenio, Alloyskull, Vylixan, Gofres, Celorico, Risque, Madgile\n\n  네오피테 티어:\n\nTacticalCrumpet, Petrie, JawnWick'

# Finish filling column B for the remaining new rows
$ws.Range('B19').Value = 'ThingDef'
$ws.Range('B20').Value = 'ThingDef'
$ws.Range('B21').Value = 'ThingDef'

# Match the author's final selection in the saved file
$ws.Range('C19').Select()
